$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wsIndi = $wb.Worksheets.Item("INDI")
Set-TextCell $wsIndi.Cells.Item(106, 1) "2022-01-01"
$wsIndi.Cells.Item(106, 2).Value = 652337.430519783
$wsIndi.Cells.Item(106, 3).Value = 65.2337430519783
Set-TextCell $wsIndi.Cells.Item(107, 1) "2022-04-01"
$wsIndi.Cells.Item(107, 2).Value = 1391602.35233014
$wsIndi.Cells.Item(107, 3).Value = 139.160235233014
Set-TextCell $wsIndi.Cells.Item(108, 1) "2022-07-01"
$wsIndi.Cells.Item(108, 2).Value = 1377764.88603498
$wsIndi.Cells.Item(108, 3).Value = 137.776488603498
Set-TextCell $wsIndi.Cells.Item(109, 1) "2022-10-01"
$wsIndi.Cells.Item(109, 2).Value = 1497787.69083112
$wsIndi.Cells.Item(109, 3).Value = 149.778769083112
Set-TextCell $wsIndi.Cells.Item(110, 1) "2023-01-01"
$wsIndi.Cells.Item(110, 2).Value = 674382.277702118
$wsIndi.Cells.Item(110, 3).Value = 67.4382277702118
Set-TextCell $wsIndi.Cells.Item(111, 1) "2023-04-01"
$wsIndi.Cells.Item(111, 2).Value = 1438670.81950154
$wsIndi.Cells.Item(111, 3).Value = 143.867081950154
Set-TextCell $wsIndi.Cells.Item(112, 1) "2023-07-01"
$wsIndi.Cells.Item(112, 2).Value = 1424390.5208906
$wsIndi.Cells.Item(112, 3).Value = 142.43905208906
Set-TextCell $wsIndi.Cells.Item(113, 1) "2023-10-01"
$wsIndi.Cells.Item(113, 2).Value = 1548491.52084353
$wsIndi.Cells.Item(113, 3).Value = 154.849152084353
Set-TextCell $wsIndi.Cells.Item(114, 1) "2024-01-01"
$wsIndi.Cells.Item(114, 2).Value = 697216.302293026
$wsIndi.Cells.Item(114, 3).Value = 69.7216302293026
Set-TextCell $wsIndi.Cells.Item(115, 1) "2024-04-01"
$wsIndi.Cells.Item(115, 2).Value = 1487388.90415494
$wsIndi.Cells.Item(115, 3).Value = 148.738890415494
Set-TextCell $wsIndi.Cells.Item(116, 1) "2024-07-01"
$wsIndi.Cells.Item(116, 2).Value = 1472628.57410196
$wsIndi.Cells.Item(116, 3).Value = 147.262857410196
Set-TextCell $wsIndi.Cells.Item(117, 1) "2024-10-01"
$wsIndi.Cells.Item(117, 2).Value = 1600934.6945855
$wsIndi.Cells.Item(117, 3).Value = 160.09346945855
Set-TextCell $wsIndi.Cells.Item(118, 1) "2025-01-01"
$wsIndi.Cells.Item(118, 2).Value = 720829.756119316
$wsIndi.Cells.Item(118, 3).Value = 72.0829756119316
Set-TextCell $wsIndi.Cells.Item(119, 1) "2025-04-01"
$wsIndi.Cells.Item(119, 2).Value = 1537764.91054443
$wsIndi.Cells.Item(119, 3).Value = 153.776491054443
Set-TextCell $wsIndi.Cells.Item(120, 1) "2025-07-01"
$wsIndi.Cells.Item(120, 2).Value = 1522505.17145261
$wsIndi.Cells.Item(120, 3).Value = 152.250517145261
Set-TextCell $wsIndi.Cells.Item(121, 1) "2025-10-01"
$wsIndi.Cells.Item(121, 2).Value = 1655157.23856895
$wsIndi.Cells.Item(121, 3).Value = 165.515723856895

$wsEta = $wb.Worksheets.Item("ETALONNAGE")
Set-TextCell $wsEta.Cells.Item(28, 1) "2022"
$wsEta.Cells.Item(28, 2).Value = 2357307.20684955
Set-TextCell $wsEta.Cells.Item(28, 3) "122.987308992901"
$wsEta.Cells.Item(28, 4).Value = 2.62687937641646
$wsEta.Cells.Item(28, 5).Value = 3.36085799827193
Set-TextCell $wsEta.Cells.Item(28, 6) "Acceptable"
Set-TextCell $wsEta.Cells.Item(29, 1) "2023"
$wsEta.Cells.Item(29, 2).Value = 2420090.74730695
Set-TextCell $wsEta.Cells.Item(29, 3) "127.148378473445"
$wsEta.Cells.Item(29, 4).Value = 2.66335844029884
$wsEta.Cells.Item(29, 5).Value = 3.38333240609754
Set-TextCell $wsEta.Cells.Item(29, 6) "Acceptable"
Set-TextCell $wsEta.Cells.Item(30, 1) "2024"
$wsEta.Cells.Item(30, 2).Value = 2485058.53192241
Set-TextCell $wsEta.Cells.Item(30, 3) "131.454211878386"
$wsEta.Cells.Item(30, 4).Value = 2.68451853252822
$wsEta.Cells.Item(30, 5).Value = 3.38646348198632
Set-TextCell $wsEta.Cells.Item(30, 6) "Acceptable"
Set-TextCell $wsEta.Cells.Item(31, 1) "2025"
$wsEta.Cells.Item(31, 2).Value = 2552234.96921028
Set-TextCell $wsEta.Cells.Item(31, 3) "135.906426917133"
$wsEta.Cells.Item(31, 4).Value = 2.70321348269804
$wsEta.Cells.Item(31, 5).Value = 3.38689417031837
Set-TextCell $wsEta.Cells.Item(31, 6) "Acceptable"

$wsPrev = $wb.Worksheets.Item("PREVISION")
Set-TextCell $wsPrev.Cells.Item(102, 1) "2022-01-01"
$wsPrev.Cells.Item(102, 2).Value = 371475.7564003
$wsPrev.Cells.Item(102, 3).Value = 240489.994115164
$wsPrev.Cells.Item(102, 4).Value = 233930.709111495
$wsPrev.Cells.Item(102, 5).Value = 137545.047288806
$wsPrev.Cells.Item(102, 6).Value = 16.3084357629946
$wsPrev.Cells.Item(102, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(102, 8).Value = 14746.3556658732
$wsPrev.Cells.Item(102, 9).Value = 3.33799595843256
$wsPrev.Cells.Item(102, 10).Value = 2.18646131641183
Set-TextCell $wsPrev.Cells.Item(102, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(103, 1) "2022-04-01"
$wsPrev.Cells.Item(103, 2).Value = 650332.366126681
$wsPrev.Cells.Item(103, 3).Value = 417440.390918196
$wsPrev.Cells.Item(103, 4).Value = 409536.043591047
$wsPrev.Cells.Item(103, 5).Value = 240796.322535634
$wsPrev.Cells.Item(103, 6).Value = 34.7900588082535
$wsPrev.Cells.Item(103, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(103, 8).Value = 11998.8411982553
$wsPrev.Cells.Item(103, 9).Value = 3.3512668379452
$wsPrev.Cells.Item(103, 10).Value = 2.68764578635323
Set-TextCell $wsPrev.Cells.Item(103, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(104, 1) "2022-07-01"
$wsPrev.Cells.Item(104, 2).Value = 645112.763756231
$wsPrev.Cells.Item(104, 3).Value = 415051.803435576
$wsPrev.Cells.Item(104, 4).Value = 406249.085390514
$wsPrev.Cells.Item(104, 5).Value = 238863.678365718
$wsPrev.Cells.Item(104, 6).Value = 34.4441221508745
$wsPrev.Cells.Item(104, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(104, 8).Value = 12050.0038182869
$wsPrev.Cells.Item(104, 9).Value = 3.36628983713656
$wsPrev.Cells.Item(104, 10).Value = 2.69425722914043
Set-TextCell $wsPrev.Cells.Item(104, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(105, 1) "2022-10-01"
$wsPrev.Cells.Item(105, 2).Value = 690386.320566341
$wsPrev.Cells.Item(105, 3).Value = 444146.629795952
$wsPrev.Cells.Item(105, 4).Value = 434759.35844633
$wsPrev.Cells.Item(105, 5).Value = 255626.962120011
$wsPrev.Cells.Item(105, 6).Value = 37.444692270778
$wsPrev.Cells.Item(105, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(105, 8).Value = 11861.4041900557
$wsPrev.Cells.Item(105, 9).Value = 3.37473497854324
$wsPrev.Cells.Item(105, 10).Value = 2.74488658045891
Set-TextCell $wsPrev.Cells.Item(105, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(106, 1) "2023-01-01"
$wsPrev.Cells.Item(106, 2).Value = 379791.248135182
$wsPrev.Cells.Item(106, 3).Value = 248617.023043927
$wsPrev.Cells.Item(106, 4).Value = 239167.252397661
$wsPrev.Cells.Item(106, 5).Value = 140623.995737521
$wsPrev.Cells.Item(106, 6).Value = 16.8595569425529
$wsPrev.Cells.Item(106, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(106, 8).Value = 14746.3556658732
$wsPrev.Cells.Item(106, 9).Value = 3.37936260452962
$wsPrev.Cells.Item(106, 10).Value = 2.23850186495644
Set-TextCell $wsPrev.Cells.Item(106, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(107, 1) "2023-04-01"
$wsPrev.Cells.Item(107, 2).Value = 668086.96639676
$wsPrev.Cells.Item(107, 3).Value = 431559.567494069
$wsPrev.Cells.Item(107, 4).Value = 420716.709245834
$wsPrev.Cells.Item(107, 5).Value = 247370.257150926
$wsPrev.Cells.Item(107, 6).Value = 35.9667704875385
$wsPrev.Cells.Item(107, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(107, 8).Value = 11998.8411982553
$wsPrev.Cells.Item(107, 9).Value = 3.38232161598366
$wsPrev.Cells.Item(107, 10).Value = 2.73008098548495
Set-TextCell $wsPrev.Cells.Item(107, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(108, 1) "2023-07-01"
$wsPrev.Cells.Item(108, 2).Value = 662700.324153847
$wsPrev.Cells.Item(108, 3).Value = 429097.780386583
$wsPrev.Cells.Item(108, 4).Value = 417324.560450377
$wsPrev.Cells.Item(108, 5).Value = 245375.76370347
$wsPrev.Cells.Item(108, 6).Value = 35.609763022265
$wsPrev.Cells.Item(108, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(108, 8).Value = 12050.0038182869
$wsPrev.Cells.Item(108, 9).Value = 3.38415032406598
$wsPrev.Cells.Item(108, 10).Value = 2.72627692176022
Set-TextCell $wsPrev.Cells.Item(108, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(109, 1) "2023-10-01"
$wsPrev.Cells.Item(109, 2).Value = 709512.208621166
$wsPrev.Cells.Item(109, 3).Value = 459182.095339979
$wsPrev.Cells.Item(109, 4).Value = 446803.569886688
$wsPrev.Cells.Item(109, 5).Value = 262708.638734477
$wsPrev.Cells.Item(109, 6).Value = 38.7122880210883
$wsPrev.Cells.Item(109, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(109, 8).Value = 11861.4041900557
$wsPrev.Cells.Item(109, 9).Value = 3.38524814450003
$wsPrev.Cells.Item(109, 10).Value = 2.77031677556052
Set-TextCell $wsPrev.Cells.Item(109, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(110, 1) "2024-01-01"
$wsPrev.Cells.Item(110, 2).Value = 388404.423866876
$wsPrev.Cells.Item(110, 3).Value = 257034.989241448
$wsPrev.Cells.Item(110, 4).Value = 244591.257253703
$wsPrev.Cells.Item(110, 5).Value = 143813.166613173
$wsPrev.Cells.Item(110, 6).Value = 17.4304075573257
$wsPrev.Cells.Item(110, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(110, 8).Value = 14746.3556658732
$wsPrev.Cells.Item(110, 9).Value = 3.38591706008533
$wsPrev.Cells.Item(110, 10).Value = 2.26787104073234
Set-TextCell $wsPrev.Cells.Item(110, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(111, 1) "2024-04-01"
$wsPrev.Cells.Item(111, 2).Value = 686463.815504575
$wsPrev.Cells.Item(111, 3).Value = 446173.581525052
$wsPrev.Cells.Item(111, 4).Value = 432289.225807032
$wsPrev.Cells.Item(111, 5).Value = 254174.589697543
$wsPrev.Cells.Item(111, 6).Value = 37.1847226038735
$wsPrev.Cells.Item(111, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(111, 8).Value = 11998.8411982553
$wsPrev.Cells.Item(111, 9).Value = 3.38632604436084
$wsPrev.Cells.Item(111, 10).Value = 2.75066720833197
Set-TextCell $wsPrev.Cells.Item(111, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(112, 1) "2024-07-01"
$wsPrev.Cells.Item(112, 2).Value = 680896.101582837
$wsPrev.Cells.Item(112, 3).Value = 443629.498521174
$wsPrev.Cells.Item(112, 4).Value = 428783.050118843
$wsPrev.Cells.Item(112, 5).Value = 252113.051463993
$wsPrev.Cells.Item(112, 6).Value = 36.815714352549
$wsPrev.Cells.Item(112, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(112, 8).Value = 12050.0038182869
$wsPrev.Cells.Item(112, 9).Value = 3.38657499498096
$wsPrev.Cells.Item(112, 10).Value = 2.74570220743777
Set-TextCell $wsPrev.Cells.Item(112, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(113, 1) "2024-10-01"
$wsPrev.Cells.Item(113, 2).Value = 729294.190968123
$wsPrev.Cells.Item(113, 3).Value = 474733.33735905
$wsPrev.Cells.Item(113, 4).Value = 459260.945848171
$wsPrev.Cells.Item(113, 5).Value = 270033.245119952
$wsPrev.Cells.Item(113, 6).Value = 40.0233673646375
$wsPrev.Cells.Item(113, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(113, 8).Value = 11861.4041900557
$wsPrev.Cells.Item(113, 9).Value = 3.38672656815078
$wsPrev.Cells.Item(113, 10).Value = 2.78811021242336
Set-TextCell $wsPrev.Cells.Item(113, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(114, 1) "2025-01-01"
$wsPrev.Cells.Item(114, 2).Value = 397311.606506902
$wsPrev.Cells.Item(114, 3).Value = 265740.298957001
$wsPrev.Cells.Item(114, 4).Value = 250200.408094012
$wsPrev.Cells.Item(114, 5).Value = 147111.198412891
$wsPrev.Cells.Item(114, 6).Value = 18.0207439029829
$wsPrev.Cells.Item(114, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(114, 8).Value = 14746.3556658732
$wsPrev.Cells.Item(114, 9).Value = 3.38681894680737
$wsPrev.Cells.Item(114, 10).Value = 2.29327528027314
Set-TextCell $wsPrev.Cells.Item(114, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(115, 1) "2025-04-01"
$wsPrev.Cells.Item(115, 2).Value = 705466.045880834
$wsPrev.Cells.Item(115, 3).Value = 461284.924046796
$wsPrev.Cells.Item(115, 4).Value = 444255.565871034
$wsPrev.Cells.Item(115, 5).Value = 261210.4800098
$wsPrev.Cells.Item(115, 6).Value = 38.4441227636107
$wsPrev.Cells.Item(115, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(115, 8).Value = 11998.8411982553
$wsPrev.Cells.Item(115, 9).Value = 3.38687523140502
$wsPrev.Cells.Item(115, 10).Value = 2.76813284940474
Set-TextCell $wsPrev.Cells.Item(115, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(116, 1) "2025-07-01"
$wsPrev.Cells.Item(116, 2).Value = 699709.950896599
$wsPrev.Cells.Item(116, 3).Value = 458654.828234136
$wsPrev.Cells.Item(116, 4).Value = 440630.760326727
$wsPrev.Cells.Item(116, 5).Value = 259079.190569872
$wsPrev.Cells.Item(116, 6).Value = 38.0626292863152
$wsPrev.Cells.Item(116, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(116, 8).Value = 12050.0038182869
$wsPrev.Cells.Item(116, 9).Value = 3.3869095186521
$wsPrev.Cells.Item(116, 10).Value = 2.76310134101614
Set-TextCell $wsPrev.Cells.Item(116, 11) "Acceptable"
Set-TextCell $wsPrev.Cells.Item(117, 1) "2025-10-01"
$wsPrev.Cells.Item(117, 2).Value = 749747.36592594
$wsPrev.Cells.Item(117, 3).Value = 490812.225119069
$wsPrev.Cells.Item(117, 4).Value = 472140.994247098
$wsPrev.Cells.Item(117, 5).Value = 277606.371678843
$wsPrev.Cells.Item(117, 6).Value = 41.3789309642237
$wsPrev.Cells.Item(117, 7).Value = 66460.5912471754
$wsPrev.Cells.Item(117, 8).Value = 11861.4041900557
$wsPrev.Cells.Item(117, 9).Value = 3.38693040801947
$wsPrev.Cells.Item(117, 10).Value = 2.80451636817045
Set-TextCell $wsPrev.Cells.Item(117, 11) "Acceptable"

$wsVatrim = $wb.Worksheets.Item("VATRIM")
Set-TextCell $wsVatrim.Cells.Item(102, 1) "2022-01-01"
$wsVatrim.Cells.Item(102, 2).Value = 240489.994115164
Set-TextCell $wsVatrim.Cells.Item(103, 1) "2022-04-01"
$wsVatrim.Cells.Item(103, 2).Value = 417440.390918196
Set-TextCell $wsVatrim.Cells.Item(104, 1) "2022-07-01"
$wsVatrim.Cells.Item(104, 2).Value = 415051.803435576
Set-TextCell $wsVatrim.Cells.Item(105, 1) "2022-10-01"
$wsVatrim.Cells.Item(105, 2).Value = 444146.629795952
Set-TextCell $wsVatrim.Cells.Item(106, 1) "2023-01-01"
$wsVatrim.Cells.Item(106, 2).Value = 248617.023043927
Set-TextCell $wsVatrim.Cells.Item(107, 1) "2023-04-01"
$wsVatrim.Cells.Item(107, 2).Value = 431559.567494069
Set-TextCell $wsVatrim.Cells.Item(108, 1) "2023-07-01"
$wsVatrim.Cells.Item(108, 2).Value = 429097.780386583
Set-TextCell $wsVatrim.Cells.Item(109, 1) "2023-10-01"
$wsVatrim.Cells.Item(109, 2).Value = 459182.095339979
Set-TextCell $wsVatrim.Cells.Item(110, 1) "2024-01-01"
$wsVatrim.Cells.Item(110, 2).Value = 257034.989241448
Set-TextCell $wsVatrim.Cells.Item(111, 1) "2024-04-01"
$wsVatrim.Cells.Item(111, 2).Value = 446173.581525052
Set-TextCell $wsVatrim.Cells.Item(112, 1) "2024-07-01"
$wsVatrim.Cells.Item(112, 2).Value = 443629.498521174
Set-TextCell $wsVatrim.Cells.Item(113, 1) "2024-10-01"
$wsVatrim.Cells.Item(113, 2).Value = 474733.33735905
Set-TextCell $wsVatrim.Cells.Item(114, 1) "2025-01-01"
$wsVatrim.Cells.Item(114, 2).Value = 265740.298957001
Set-TextCell $wsVatrim.Cells.Item(115, 1) "2025-04-01"
$wsVatrim.Cells.Item(115, 2).Value = 461284.924046796
Set-TextCell $wsVatrim.Cells.Item(116, 1) "2025-07-01"
$wsVatrim.Cells.Item(116, 2).Value = 458654.828234136
Set-TextCell $wsVatrim.Cells.Item(117, 1) "2025-10-01"
$wsVatrim.Cells.Item(117, 2).Value = 490812.225119069

Write-Host "Edit complete"